$wb = $excel.ActiveWorkbook

# --- 1. Rename the "Include from ParticipationSig" sheet to "Include #0" ---
$wsInclude = $wb.Worksheets.Item(2)
$wsInclude.Name = "Include #0"

# --- 2. Update the Metadata sheet ---
$wsMeta = $wb.Worksheets.Item(1)

# Version value (row 3)
$wsMeta.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"

# Date value (row 8)
$wsMeta.Range("B8").Value = "2025-10-29T22:15:57+01:00"

# Insert a new row after "Contact" (row 10) for "Jurisdiction", pushing
# Description/Purpose/Copyright/Immutable down by one row.
$wsMeta.Rows.Item(11).Insert()

# Copy the bordered style from the row that landed below (old row 11, now
# row 12) onto the freshly inserted blank row so it matches the rest of
# the table instead of keeping Excel's default insert style.
$wsMeta.Range("A12:B12").Copy()
$wsMeta.Range("A11:B11").PasteSpecial(-4122)

# Populate the new "Jurisdiction" row (value left blank).
$wsMeta.Range("A11").Value = "Jurisdiction"
$wsMeta.Range("B11").Value = ""

Write-Host "done"
